$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from row 5 to new row 6 so the new entry matches existing style
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122) # xlPasteFormats

$ws.Range("B5:H5").Copy()
$ws.Range("B6:H6").PasteSpecial(-4122) # xlPasteFormats

# Populate new row 6 values in the order the shared strings were added
$ws.Range("A6").Value = 40918
$ws.Range("C6").Value = "1 phần SDD"
$ws.Range("D6").Value = "chưa hoàn thành"
$ws.Range("B6").Value = "viết các chức năng SDD như SRS"
$ws.Range("E6").Value = 4

$ws.Range("H9").Select()
